$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure D and E columns are stored as text so numeric-looking values
# (e.g. "1.0000", "239.04") keep their exact original formatting instead
# of being auto-converted to numbers by Excel.
$ws.Range("D2:E52").NumberFormat = "@"

# --- Refresh Price / Volume(1h) for the coins that stayed on rows 2-43 ---
# --- Coin list shifted down by one row starting at row 44 (new entry inserted) ---
$ws.Range("D2").Value = "25.821.40"
$ws.Range("E2").Value = "  -2.49%  "
$ws.Range("D3").Value = "1.750.71"
$ws.Range("E3").Value = "  -4.67%  "
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "239.04"
$ws.Range("E5").Value = "  -8.45%  "
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "0.5090"
$ws.Range("E7").Value = "  -5.30%  "
$ws.Range("D8").Value = "42.10"
$ws.Range("E8").Value = "  -6.16%  "
$ws.Range("D9").Value = "0.2769"
$ws.Range("E9").Value = "  -6.64%  "
$ws.Range("D10").Value = "0.06191"
$ws.Range("E10").Value = "  -10.46%  "
$ws.Range("D11").Value = "1.747.75"
$ws.Range("E11").Value = "  -5.33%  "
$ws.Range("D12").Value = "0.06958"
$ws.Range("E12").Value = "  -3.49%  "
$ws.Range("D13").Value = "15.72"
$ws.Range("E13").Value = "  -10.32%  "
$ws.Range("D14").Value = "0.6051"
$ws.Range("E14").Value = "  -17.36%  "
$ws.Range("D15").Value = "4.517"
$ws.Range("E15").Value = "  -9.44%  "
$ws.Range("D16").Value = "77.63"
$ws.Range("E16").Value = "  -12.80%  "
$ws.Range("D17").Value = "0.9998"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "0.9993"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "25.830.65"
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("D20").Value = "0.000006966"
$ws.Range("E20").Value = "  -11.72%  "
$ws.Range("D21").Value = "11.68"
$ws.Range("E21").Value = "  -15.43%  "
$ws.Range("D22").Value = "1.970.31"
$ws.Range("E22").Value = "  -5.16%  "
$ws.Range("D23").Value = "4.087"
$ws.Range("E23").Value = "  -10.68%  "
$ws.Range("D24").Value = "5.245"
$ws.Range("E24").Value = "  -12.28%  "
$ws.Range("D25").Value = "8.217"
$ws.Range("E25").Value = "  -10.62%  "
$ws.Range("D26").Value = "137.66"
$ws.Range("E26").Value = "  -3.34%  "
$ws.Range("D27").Value = "1.473"
$ws.Range("E27").Value = "  -13.76%  "
$ws.Range("D28").Value = "1.823"
$ws.Range("E28").Value = "  -15.92%  "
$ws.Range("E29").Value = "  -11.35%  "
$ws.Range("D30").Value = "103.77"
$ws.Range("E30").Value = "  -6.35%  "
$ws.Range("D31").Value = "0.08166"
$ws.Range("E31").Value = "  -7.71%  "
$ws.Range("D32").Value = "3.703"
$ws.Range("E32").Value = "  -12.35%  "
$ws.Range("D33").Value = "3.516"
$ws.Range("E33").Value = "  -12.83%  "
$ws.Range("D34").Value = "0.04536"
$ws.Range("E34").Value = "  -6.19%  "
$ws.Range("D35").Value = "0.9984"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "2.625"
$ws.Range("E36").Value = "  -10.40%  "
$ws.Range("D37").Value = "0.9901"
$ws.Range("E37").Value = "  -12.34%  "
$ws.Range("D38").Value = "0.6117"
$ws.Range("E38").Value = "  -15.47%  "
$ws.Range("D39").Value = "2.695"
$ws.Range("E39").Value = "  -12.77%  "
$ws.Range("D40").Value = "0.01557"
$ws.Range("E40").Value = "  -9.07%  "
$ws.Range("D41").Value = "1.914"
$ws.Range("E41").Value = "  -16.32%  "
$ws.Range("D42").Value = "0.9992"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "103.10"
$ws.Range("E43").Value = "  -4.15%  "
$ws.Range("B44").Value = "PaxosStandard"
$ws.Range("C44").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.3879"
$ws.Range("E45").Value = "  -17.40%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "0.7474"
$ws.Range("E46").Value = "  -17.30%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "4.935"
$ws.Range("E47").Value = "  -16.04%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05398"
$ws.Range("E48").Value = "  -6.35%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.1114"
$ws.Range("E49").Value = "  -10.65%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "6.009"
$ws.Range("E50").Value = "  -18.69%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "30.18"
$ws.Range("E51").Value = "  -13.12%  "

# --- New row 52: the last coin (Aave) that got pushed out of row 51 ---
$ws.Range("A51").Copy()
$ws.Range("A52").PasteSpecial(-4122)  # xlPasteFormats, match column-A styling
$ws.Range("A52").Value = 50
$ws.Range("B52").Value = "Aave"
$ws.Range("C52").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D52").Value = "52.61"
$ws.Range("E52").Value = "  -12.30%  "
